# Weekly price-sheet update: a new week's worth of data (2 rows) is
# inserted at the top of the "Primera" date-block starting at row 342,
# pushing the existing rows 342:435 down to 344:437.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 342/343; Excel shifts 342:435 -> 344:437
# and grows the used range to A1:R437 automatically.
$ws.Range("A342:A343").EntireRow.Insert()

# ---- New row 342 ----
$ws.Cells.Item(342, 1).Value = 9
$ws.Cells.Item(342, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(342, 3).Value = "Metropolitana"
$ws.Cells.Item(342, 4).Value = 44551
$ws.Cells.Item(342, 5).Value = 13
$ws.Cells.Item(342, 6).Value = 100112040
$ws.Cells.Item(342, 7).Value = "Cilantro"
$ws.Cells.Item(342, 8).Value = "Sin especificar"
$ws.Cells.Item(342, 9).Value = "Primera"
$ws.Cells.Item(342, 10).Value = 43
$ws.Cells.Item(342, 11).Value = 8000
$ws.Cells.Item(342, 12).Value = 8000
$ws.Cells.Item(342, 13).Value = 8000
$ws.Cells.Item(342, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(342, 15).Value = "Región Metropolitana"
$ws.Cells.Item(342, 16).Value = 222
$ws.Cells.Item(342, 17).Value = 36
$ws.Cells.Item(342, 18).Value = "Hortaliza"

# ---- New row 343 ----
$ws.Cells.Item(343, 1).Value = 9
$ws.Cells.Item(343, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(343, 3).Value = "Metropolitana"
$ws.Cells.Item(343, 4).Value = 44551
$ws.Cells.Item(343, 5).Value = 13
$ws.Cells.Item(343, 6).Value = 100112040
$ws.Cells.Item(343, 7).Value = "Cilantro"
$ws.Cells.Item(343, 8).Value = "Sin especificar"
$ws.Cells.Item(343, 9).Value = "Primera"
$ws.Cells.Item(343, 10).Value = 106
$ws.Cells.Item(343, 11).Value = 17000
$ws.Cells.Item(343, 12).Value = 19000
$ws.Cells.Item(343, 13).Value = 18000
$ws.Cells.Item(343, 14).Value = "$/docena de atados"
$ws.Cells.Item(343, 15).Value = "Región Metropolitana"
$ws.Cells.Item(343, 16).Value = 6000
$ws.Cells.Item(343, 17).Value = 3
$ws.Cells.Item(343, 18).Value = "Hortaliza"
